# Generate Report for Handback
# The file "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md" has now been handed back
# (in both zh-cn and de-de). The status tables are re-sorted alphabetically by
# source file name (since "13d13c86..." < "ffffef..." < "ffffff..." in ASCII),
# and the 13d13c86 row picks up its "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" data, while the "Ready for handoff"
# status throughout is replaced by "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
$ov.Range("A3").Value = "ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md"
$ov.Range("A4").Value = "ffffff6346cb71-d9d5-42e5-be54-4faf74872423.md"

$ov.Range("B4").Value = "Handed back: in sync with en-US"
$ov.Range("C4").Value = "Handed back: in sync with en-US"

foreach ($hl in $ov.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$2") {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/54e9cae050ae73c8d3a3c3adf832ebdfd4459756/e2e/13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
        $hl.TextToDisplay = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
    } elseif ($addr -eq "`$A`$3") {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/54e9cae050ae73c8d3a3c3adf832ebdfd4459756/e2e/ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md"
        $hl.TextToDisplay = "ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md"
    } elseif ($addr -eq "`$A`$4") {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/54e9cae050ae73c8d3a3c3adf832ebdfd4459756/e2e/ffffff6346cb71-d9d5-42e5-be54-4faf74872423.md"
        $hl.TextToDisplay = "ffffff6346cb71-d9d5-42e5-be54-4faf74872423.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
$zh.Range("C2").Value = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.zh-cn.xlf"
$zh.Range("D2").Value = "2016-01-25 08:26:40"
$zh.Range("E2").Value = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
$zh.Range("F2").Value = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.zh-cn.xlf"
$zh.Range("G2").Value = "2016-01-25 08:27:24"

$zh.Range("A3").Value = "ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md"

$zh.Range("A4").Value = "ffffff6346cb71-d9d5-42e5-be54-4faf74872423.md"
$zh.Range("B4").Value = "Handed back: in sync with en-US"
$zh.Range("C4").Value = "fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.zh-cn.xlf"
$zh.Range("D4").Value = "2016-01-25 08:22:33"
$zh.Range("E4").Value = "fc1a2693-ec42-4cdf-a862-a22868795b5c.md"
$zh.Range("F4").Value = "fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.zh-cn.xlf"
$zh.Range("G4").Value = "2016-01-25 08:23:17"

foreach ($hl in $zh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$2") {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/54e9cae050ae73c8d3a3c3adf832ebdfd4459756/e2e/13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
        $hl.TextToDisplay = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
    } elseif ($addr -eq "`$C`$2") {
        $hl.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4b69e0b13467f024a497858b7960bfbf92dac33a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.zh-cn.xlf"
        $hl.TextToDisplay = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.zh-cn.xlf"
    } elseif ($addr -eq "`$E`$2") {
        $hl.Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c2200d014a5873c2c96c4ceec9213ec1ad02fc9c/e2e/13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
        $hl.TextToDisplay = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
    } elseif ($addr -eq "`$F`$2") {
        $hl.Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/553ec1072e5850251f7b04db472ece1d1aa5a096/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.zh-cn.xlf"
        $hl.TextToDisplay = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.zh-cn.xlf"
    } elseif ($addr -eq "`$A`$3") {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/54e9cae050ae73c8d3a3c3adf832ebdfd4459756/e2e/ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md"
        $hl.TextToDisplay = "ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md"
    } elseif ($addr -eq "`$A`$4") {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/54e9cae050ae73c8d3a3c3adf832ebdfd4459756/e2e/ffffff6346cb71-d9d5-42e5-be54-4faf74872423.md"
        $hl.TextToDisplay = "ffffff6346cb71-d9d5-42e5-be54-4faf74872423.md"
    } elseif ($addr -eq "`$C`$4") {
        $hl.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34d41e7f2edce79281388abbdbd7bdb92bb98ae6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.zh-cn.xlf"
        $hl.TextToDisplay = "fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.zh-cn.xlf"
    }
}

# New hyperlinks for E4 / F4 (these cells had no hyperlink before the handback)
$zh.Hyperlinks.Add($zh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c2200d014a5873c2c96c4ceec9213ec1ad02fc9c/e2e/fc1a2693-ec42-4cdf-a862-a22868795b5c.md", "", "", "fc1a2693-ec42-4cdf-a862-a22868795b5c.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/553ec1072e5850251f7b04db472ece1d1aa5a096/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.zh-cn.xlf", "", "", "fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.zh-cn.xlf") | Out-Null
$zh.Range("E4").Style = $zh.Range("E3").Style
$zh.Range("F4").Style = $zh.Range("F3").Style

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
$de.Range("C2").Value = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.de-de.xlf"
$de.Range("D2").Value = "2016-01-25 08:26:52"
$de.Range("E2").Value = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
$de.Range("F2").Value = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.de-de.xlf"
$de.Range("G2").Value = "2016-01-25 08:27:46"

$de.Range("A3").Value = "ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md"

$de.Range("A4").Value = "ffffff6346cb71-d9d5-42e5-be54-4faf74872423.md"
$de.Range("B4").Value = "Handed back: in sync with en-US"
$de.Range("C4").Value = "fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.de-de.xlf"
$de.Range("D4").Value = "2016-01-25 08:22:45"
$de.Range("E4").Value = "fc1a2693-ec42-4cdf-a862-a22868795b5c.md"
$de.Range("F4").Value = "fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.de-de.xlf"
$de.Range("G4").Value = "2016-01-25 08:23:39"

foreach ($hl in $de.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$2") {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/54e9cae050ae73c8d3a3c3adf832ebdfd4459756/e2e/13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
        $hl.TextToDisplay = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
    } elseif ($addr -eq "`$C`$2") {
        $hl.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2031e4cda1d13ea739ac175b650fe6fa8f370080/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.de-de.xlf"
        $hl.TextToDisplay = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.de-de.xlf"
    } elseif ($addr -eq "`$E`$2") {
        $hl.Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/83e7b41a8f3223463026f945b77e4f8e284f9fd1/e2e/13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
        $hl.TextToDisplay = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.md"
    } elseif ($addr -eq "`$F`$2") {
        $hl.Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9def0246a9a3bfe2f16e222c175ac24fd6c44789/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.de-de.xlf"
        $hl.TextToDisplay = "13d13c86-4fdf-4468-bfd0-b3bbae73354e.05bf5bd1cf312a1cfee430b5de49093266ba27fd.de-de.xlf"
    } elseif ($addr -eq "`$A`$3") {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/54e9cae050ae73c8d3a3c3adf832ebdfd4459756/e2e/ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md"
        $hl.TextToDisplay = "ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md"
    } elseif ($addr -eq "`$A`$4") {
        $hl.Address = "https://github.com/OpenLocalizationTest/oltest/blob/54e9cae050ae73c8d3a3c3adf832ebdfd4459756/e2e/ffffff6346cb71-d9d5-42e5-be54-4faf74872423.md"
        $hl.TextToDisplay = "ffffff6346cb71-d9d5-42e5-be54-4faf74872423.md"
    } elseif ($addr -eq "`$C`$4") {
        $hl.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82ec11bc4fbb07fcd56c85145219ad67af0c9d1b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.de-de.xlf"
        $hl.TextToDisplay = "fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.de-de.xlf"
    }
}

# New hyperlinks for E4 / F4 (these cells had no hyperlink before the handback)
$de.Hyperlinks.Add($de.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/83e7b41a8f3223463026f945b77e4f8e284f9fd1/e2e/fc1a2693-ec42-4cdf-a862-a22868795b5c.md", "", "", "fc1a2693-ec42-4cdf-a862-a22868795b5c.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9def0246a9a3bfe2f16e222c175ac24fd6c44789/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.de-de.xlf", "", "", "fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5.de-de.xlf") | Out-Null
$de.Range("E4").Style = $de.Range("E3").Style
$de.Range("F4").Style = $de.Range("F3").Style
